$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$labels = @("model_10_9_0", "model_10_9_1", "model_10_9_2", "model_10_9_3", "model_10_9_4", "model_10_9_5", "model_10_9_6", "model_10_9_7", "model_10_9_11", "model_10_9_22", "model_10_9_23", "model_10_9_24", "model_10_9_14", "model_10_9_12", "model_10_9_13", "model_10_9_18", "model_10_9_19", "model_10_9_15", "model_10_9_17", "model_10_9_20", "model_10_9_10", "model_10_9_9", "model_10_9_8", "model_10_9_16", "model_10_9_21")
$Bvals = @(0.6731329872420542, 0.8136900225301507, 0.8828588124210566, 0.9052290444307735, 0.9101162290291935, 0.9142427645492304, 0.9280227511191278, 0.9337554629182967, 0.9359938362635287, 0.9371250915051701, 0.9372992084207223, 0.937394534616672, 0.9374879410103973, 0.937511017889053, 0.9376725939094953, 0.9377787102141649, 0.9378472007313067, 0.9381102395834383, 0.9381741652246556, 0.938177415345669, 0.9385606148622494, 0.9385837610998704, 0.9388660547031141, 0.9390298905720929, 0.9403143444439802)
$Cvals = @(0.9642779591401296, 0.8460423068903125, 0.9221932278272831, 0.9242002865580485, 0.9243855694987653, 0.8894638334525801, 0.9073425709233736, 0.9230164725916787, 0.9291807427249357, 0.9231506113735051, 0.9231389275941588, 0.9234029459860633, 0.9341816456737193, 0.9360555593190121, 0.9365106030120325, 0.9309741351546748, 0.9303846036080728, 0.9350288837742073, 0.9332294909965185, 0.9315489676816726, 0.941694923576046, 0.9423257850840271, 0.9422497585390088, 0.9385558291782997, 0.9380783388551716)
$Dvals = @(0.9902500836463899, 0.9650006768277025, 0.9383688514485632, 0.9400081740137053, 0.9423666186110607, 0.9418100035385419, 0.9406227109147, 0.9375792447538449, 0.91660063375879, 0.9240739102826842, 0.9250119012714355, 0.9251750076542413, 0.9155785644568722, 0.9149501353479504, 0.9143707358146593, 0.9189640572159539, 0.9195966080904711, 0.9180690761916371, 0.9188703910562553, 0.9189931654486759, 0.9240305949210121, 0.9240968985131502, 0.9326862738462551, 0.9180812375328996, 0.923628156358609)
$Evals = @(0.9738501082715144, 0.8898459338497586, 0.9283805755560042, 0.9302460583351643, 0.9312190184974655, 0.9089123457188047, 0.9197980917795423, 0.9286143699434306, 0.9249022875045758, 0.9237976430731787, 0.9241312148279113, 0.9243573539804931, 0.9276968384486661, 0.9286545723003305, 0.9287322711918287, 0.9268952738814328, 0.9267510435523931, 0.9291367251454857, 0.9282889468723197, 0.92727000389701, 0.9355191770179158, 0.9359434361580654, 0.9390128075158672, 0.9313732255364282, 0.9330842910288129)
$Fvals = @(0.3617455065250397, 0.2061902433633804, 0.1296407729387283, 0.1048835292458534, 0.09947484731674194, 0.09490799158811569, 0.07965760678052902, 0.07331318408250809, 0.07083596289157867, 0.06958399713039398, 0.06939129531383514, 0.0692858025431633, 0.06918243318796158, 0.06915689259767532, 0.06897807121276855, 0.06886063516139984, 0.06878483295440674, 0.06849373131990433, 0.06842298060655594, 0.06841938197612762, 0.06799529492855072, 0.0679696798324585, 0.06765726208686829, 0.06747595220804214, 0.06605443358421326)
$Gvals = @(0.05636938288807869, 0.2429452538490295, 0.122779093682766, 0.1196119487285614, 0.1193195730447769, 0.1744260936975479, 0.1462134271860123, 0.1214800179004669, 0.1117528006434441, 0.1212683394551277, 0.1212867796421051, 0.1208701506257057, 0.1038613691926003, 0.1009043306112289, 0.1001862734556198, 0.1089228242635727, 0.1098531037569046, 0.1025244370102882, 0.1053638756275177, 0.1080157309770584, 0.09200542420148849, 0.09100992232561111, 0.09112988412380219, 0.09695890545845032, 0.09771239012479782)
$Hvals = @(0.009925415739417076, 0.03562930971384048, 0.0627405121922493, 0.06107167899608612, 0.0586707815527916, 0.05923741310834885, 0.06044608354568481, 0.06354433298110962, 0.08490055799484253, 0.07729276269674301, 0.07633788883686066, 0.07617184519767761, 0.08594102412462234, 0.08658076077699661, 0.08717059344053268, 0.0824945867061615, 0.08185065537691116, 0.08340568840503693, 0.08258993923664093, 0.08246495574712753, 0.07733685523271561, 0.07726936042308807, 0.06852537393569946, 0.08339329808950424, 0.07774654030799866)
$Ivals = @(0.03451340645551682, 0.1453846246004105, 0.09452544897794724, 0.09206333011388779, 0.09077919274568558, 0.1202201992273331, 0.105852872133255, 0.09421688318252563, 0.09911619871854782, 0.1005741506814957, 0.1001338809728622, 0.09983541816473007, 0.09542787075042725, 0.09416382759809494, 0.09406127035617828, 0.09648579359054565, 0.09667615592479706, 0.09352746605873108, 0.0946463868021965, 0.09599120914936066, 0.08510371297597885, 0.08454376459121704, 0.08049271255731583, 0.09057566523551941, 0.08831734955310822)

for ($i = 0; $i -lt 25; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $labels[$i]
    $ws.Cells.Item($row, 2).Value = $Bvals[$i]
    $ws.Cells.Item($row, 3).Value = $Cvals[$i]
    $ws.Cells.Item($row, 4).Value = $Dvals[$i]
    $ws.Cells.Item($row, 5).Value = $Evals[$i]
    $ws.Cells.Item($row, 6).Value = $Fvals[$i]
    $ws.Cells.Item($row, 7).Value = $Gvals[$i]
    $ws.Cells.Item($row, 8).Value = $Hvals[$i]
    $ws.Cells.Item($row, 9).Value = $Ivals[$i]
}
